# Apply updated odds values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.47
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 2.75
$ws.Range("J2").Value = 3.05
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 3.4
$ws.Range("M2").Value = 9.4
$ws.Range("N2").Value = 1.03
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 1.91
$ws.Range("R2").Value = 1.8
$ws.Range("S2").Value = 1.42
$ws.Range("T2").Value = 2.47
$ws.Range("U2").Value = 1.65
$ws.Range("V2").Value = 1.98
$ws.Range("W2").Value = 8.5
$ws.Range("X2").Value = 13
$ws.Range("Y2").Value = 9.25
$ws.Range("Z2").Value = 27
$ws.Range("AA2").Value = 19.5
$ws.Range("AB2").Value = 27
$ws.Range("AC2").Value = 9.5
$ws.Range("AE2").Value = 13
$ws.Range("AH2").Value = 8.75
$ws.Range("AI2").Value = 14
$ws.Range("AK2").Value = 32
$ws.Range("AN2").Value = 4.35
$ws.Range("AP2").Value = 21
$ws.Range("AT2").Value = 2.45
$ws.Range("AU2").Value = 6.8
$ws.Range("AW2").Value = 4.65

# Row 3
$ws.Range("G3").Value = 2.05
$ws.Range("H3").Value = 2.95
$ws.Range("I3").Value = 3.75
$ws.Range("J3").Value = 2.6
$ws.Range("K3").Value = 2.02
$ws.Range("L3").Value = 4.15
$ws.Range("O3").Value = 1.35
$ws.Range("T3").Value = 2.52
$ws.Range("U3").Value = 1.78
$ws.Range("V3").Value = 1.83
$ws.Range("W3").Value = 6.6
$ws.Range("X3").Value = 9.5
$ws.Range("Z3").Value = 19.5
$ws.Range("AA3").Value = 17.5
$ws.Range("AB3").Value = 30
$ws.Range("AC3").Value = 7.8
$ws.Range("AD3").Value = 5.8
$ws.Range("AH3").Value = 9.75
$ws.Range("AI3").Value = 21
$ws.Range("AJ3").Value = 12.5
$ws.Range("AK3").Value = 60
$ws.Range("AN3").Value = 3.9
$ws.Range("AO3").Value = 10.5
$ws.Range("AP3").Value = 18
$ws.Range("AT3").Value = 2.5
$ws.Range("AW3").Value = 5.6
$ws.Range("AX3").Value = 21
$ws.Range("AY3").Value = 25
$ws.Range("BA3").Value = 120

# Row 7
$ws.Range("G7").Value = 3.6
$ws.Range("I7").Value = 2.05
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.75
$ws.Range("AC7").Value = 7.5
$ws.Range("AF7").Value = 67
$ws.Range("AI7").Value = 9
$ws.Range("AQ7").Value = 81
$ws.Range("AS7").Value = 301

# Row 9
$ws.Range("Q9").Value = 1.98
$ws.Range("R9").Value = 1.88

# Row 11
$ws.Range("G11").Value = 8.5
$ws.Range("H11").Value = 4.9
$ws.Range("J11").Value = 6.6
$ws.Range("K11").Value = 2.67
$ws.Range("R11").Value = 2.87
$ws.Range("S11").Value = 1.22
$ws.Range("T11").Value = 3.85
$ws.Range("U11").Value = 1.6
$ws.Range("V11").Value = 2.22
$ws.Range("W11").Value = 35
$ws.Range("X11").Value = 80
$ws.Range("Y11").Value = 25
$ws.Range("Z11").Value = 250
$ws.Range("AA11").Value = 80
$ws.Range("AD11").Value = 10.75
$ws.Range("AE11").Value = 15.5
$ws.Range("AG11").Value = 250
$ws.Range("AH11").Value = 11
$ws.Range("AJ11").Value = 8.5
$ws.Range("AM11").Value = 18
$ws.Range("AN11").Value = 10
$ws.Range("AO11").Value = 40
$ws.Range("AP11").Value = 30
$ws.Range("AQ11").Value = 250
$ws.Range("AR11").Value = 175
$ws.Range("AT11").Value = 3.85
$ws.Range("AU11").Value = 7.1
$ws.Range("AV11").Value = 45
$ws.Range("AW11").Value = 3.6
$ws.Range("AX11").Value = 5.9
$ws.Range("AZ11").Value = 14
$ws.Range("BA11").Value = 28
